$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18; Excel shifts existing rows 18..96 down to 19..97
# and copies the row-above formatting (so D18 inherits the date number format).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record.
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = "Macroferia Regional de Talca"
$ws.Range("C18").Value = "Maule"
$ws.Range("D18").Value = 44547
$ws.Range("E18").Value = 7
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100108
$ws.Range("H18").Value = "Tropicales y subtropicales"
$ws.Range("I18").Value = 100108002
$ws.Range("J18").Value = "Mango"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 210
$ws.Range("N18").Value = 7000
$ws.Range("O18").Value = 7000
$ws.Range("P18").Value = 7000
$ws.Range("Q18").Value = '$/bandeja 4 kilos'
$ws.Range("R18").Value = "Perú"
$ws.Range("S18").Value = 1750
$ws.Range("T18").Value = 4
